$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 174.88889
$ws.Range("I6").Value = 174.88889
$ws.Range("K6").Value = 524.6666700000001
$ws.Range("M6").Value = -412.6666700000001

$ws.Range("H31").Value = 799.5
$ws.Range("I31").Value = 799.5
$ws.Range("K31").Value = 2398.5
$ws.Range("M31").Value = -2168.5

$ws.Range("H40").Value = 3700
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 3800
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 3800
$ws.Range("M40").Value = -3425
$ws.Range("N40").Value = -4150

$ws.Range("H61").Value = 1000
$ws.Range("J61").Value = 1000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3344

$ws.Range("H74").Value = 4003
$ws.Range("I74").Value = 4003
$ws.Range("K74").Value = 4003
$ws.Range("M74").Value = -3067

$ws.Range("H77").Value = 4003
$ws.Range("I77").Value = 4003
$ws.Range("K77").Value = 20015
$ws.Range("M77").Value = -15335

$ws.Range("H98").Value = 1332.8462
$ws.Range("I98").Value = 1332.8462
$ws.Range("K98").Value = 1332.8462
$ws.Range("M98").Value = 165.1538

$ws.Range("H107").Value = 75552.5
$ws.Range("I107").Value = 100382.11
$ws.Range("K107").Value = 100382.11
$ws.Range("M107").Value = -98462.11

$ws.Range("H116").Value = 2987.5
$ws.Range("I116").Value = 2985
$ws.Range("J116").Value = 2990
$ws.Range("K116").Value = 2985
$ws.Range("L116").Value = 2990
$ws.Range("M116").Value = 457
$ws.Range("N116").Value = -9874

$ws.Range("H120").Value = 59989
$ws.Range("J120").Value = 59989
$ws.Range("L120").Value = 59989
$ws.Range("N120").Value = -69665

$ws.Range("H122").Value = 1332.8462
$ws.Range("I122").Value = 1332.8462
$ws.Range("K122").Value = 3998.5386
$ws.Range("M122").Value = -1548.5386

$ws.Range("H132").Value = 2042.1538
$ws.Range("I132").Value = 2042.1538
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6126.4614
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3596.4614
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 750
$ws.Range("K3").Value = 750
$ws.Range("M3").Value = -635

$ws.Range("H32").Value = 12333.333
$ws.Range("I32").Value = 12333.333
$ws.Range("K32").Value = 12333.333
$ws.Range("M32").Value = -12046.333

$ws.Range("H97").Value = 66670230
$ws.Range("I97").Value = 66670230
$ws.Range("K97").Value = 66670230
$ws.Range("M97").Value = -66669734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H99").Value = 2175.6667
$ws.Range("I99").Value = 2175.6667
$ws.Range("K99").Value = 2175.6667
$ws.Range("M99").Value = -677.6667000000002

$ws.Range("H134").Value = 3773.75
$ws.Range("I134").Value = 3773.75
$ws.Range("K134").Value = 11321.25
$ws.Range("M134").Value = -8786.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 599.25
$ws.Range("I22").Value = 603.1429000000001
$ws.Range("J22").Value = 572
$ws.Range("K22").Value = 603.1429000000001
$ws.Range("L22").Value = 572
$ws.Range("M22").Value = -253.1429000000001
$ws.Range("N22").Value = -1272

$ws.Range("H31").Value = 7181.8184
$ws.Range("I31").Value = 4002.75
$ws.Range("J31").Value = 8998.429
$ws.Range("K31").Value = 4002.75
$ws.Range("L31").Value = 8998.429
$ws.Range("M31").Value = -3707.75
$ws.Range("N31").Value = -9588.429

$ws.Range("H34").Value = 7181.8184
$ws.Range("I34").Value = 4002.75
$ws.Range("J34").Value = 8998.429
$ws.Range("K34").Value = 4002.75
$ws.Range("L34").Value = 8998.429
$ws.Range("M34").Value = -3800.75
$ws.Range("N34").Value = -9402.429

$ws.Range("H99").Value = 2472281
$ws.Range("I99").Value = 1131128.4
$ws.Range("J99").Value = 3142857.2
$ws.Range("K99").Value = 1131128.4
$ws.Range("L99").Value = 3142857.2
$ws.Range("M99").Value = -1129630.4
$ws.Range("N99").Value = -3145853.2

$ws.Range("H126").Value = 2472281
$ws.Range("I126").Value = 1131128.4
$ws.Range("J126").Value = 3142857.2
$ws.Range("K126").Value = 3393385.2
$ws.Range("L126").Value = 9428571.600000001
$ws.Range("M126").Value = -3390915.2
$ws.Range("N126").Value = -9433511.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 801.44446
$ws.Range("I6").Value = 42.6
$ws.Range("J6").Value = 1750
$ws.Range("K6").Value = 127.8
$ws.Range("L6").Value = 5250
$ws.Range("M6").Value = -14.80000000000001
$ws.Range("N6").Value = -5476

$ws.Range("H17").Value = 604.1
$ws.Range("I17").Value = 699.6667
$ws.Range("J17").Value = 563.1429000000001
$ws.Range("K17").Value = 2099.0001
$ws.Range("L17").Value = 1689.4287
$ws.Range("M17").Value = -1930.0001
$ws.Range("N17").Value = -2027.4287

$ws.Range("H26").Value = 690.2
$ws.Range("I26").Value = 690.2
$ws.Range("K26").Value = 2070.6
$ws.Range("M26").Value = -1782.6

$ws.Range("H99").Value = 4077.5715
$ws.Range("I99").Value = 4077.5715
$ws.Range("K99").Value = 12232.7145
$ws.Range("M99").Value = -9986.7145

$ws.Range("H108").Value = 427
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H131").Value = 995
$ws.Range("I131").Value = 988.6667
$ws.Range("J131").Value = 1033
$ws.Range("K131").Value = 2966.0001
$ws.Range("L131").Value = 3099
$ws.Range("M131").Value = 2073.9999
$ws.Range("N131").Value = -13179

$ws.Range("H134").Value = 875
$ws.Range("I134").Value = 875
$ws.Range("K134").Value = 2625
$ws.Range("M134").Value = 2445

$ws.Range("H139").Value = 4072.5557
$ws.Range("I139").Value = 4072.5557
$ws.Range("K139").Value = 12217.6671
$ws.Range("M139").Value = -7077.667099999999

$ws.Range("H140").Value = 669.7143
$ws.Range("I140").Value = 448
$ws.Range("K140").Value = 1344
$ws.Range("M140").Value = 3836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 29750
$ws.Range("J95").Value = 29750
$ws.Range("L95").Value = 29750
$ws.Range("N95").Value = -35242

$ws.Range("H97").Value = 3077.4
$ws.Range("I97").Value = 3575
$ws.Range("K97").Value = 3575
$ws.Range("M97").Value = -3079

$ws.Range("H114").Value = 5000
$ws.Range("J114").Value = 5000
$ws.Range("L114").Value = 5000
$ws.Range("N114").Value = -13678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 17559
$ws.Range("J106").Value = 17559
$ws.Range("L106").Value = 17559
$ws.Range("N106").Value = -20083

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 23166.334
$ws.Range("J51").Value = 22249.5
$ws.Range("L51").Value = 22249.5
$ws.Range("N51").Value = -23269.5

$ws.Range("H52").Value = 24999.5
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 39999
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 39999
$ws.Range("M52").Value = -9774
$ws.Range("N52").Value = -40451

$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178

$ws.Range("H132").Value = 1502
$ws.Range("I132").Value = 1004
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3012
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -482
$ws.Range("N132").Value = -11060
